# Generate Report for Handoff
#
# Replaces the two "in progress" source-file rows (keyed by UUID-named
# markdown files) with a new pair of rows describing the file that is now
# actually ready for/at handoff, across all three sheets (Overview, zh-cn,
# de-de), and drops the now-unused "Latest Handoff File" / "Latest Handback
# File" columns (F/G) from the per-locale sheets.

$wb = $excel.ActiveWorkbook

# New identifiers that replace the old ones everywhere.
$newFile1 = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.md"
$newFile2 = "ffffaf1f2894-e07f-4587-80c2-b4391c75d2cf.md"
$newStatus = "Ready for handoff"
$newHandoffDate = "2016-43-11 16:43:14"

$newZhTarget = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.zh-cn.xlf"
$newDeTarget = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.de-de.xlf"
$newHandoffDatetime = "2016-03-11 16:43:11"
$newHandbackDatetime = "0001-01-01 00:00:00"
$newDeHandoffDatetime = "2016-03-11 16:43:14"

$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/0a9d409a4894333fd7f8c33dea6d8204c9d1c691/e2e/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.md"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/0a9d409a4894333fd7f8c33dea6d8204c9d1c691/e2e/ffffaf1f2894-e07f-4587-80c2-b4391c75d2cf.md"
$zhTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22b0e15878f591fbd7e09c929416f6243269feee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.zh-cn.xlf"
$deTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee6d6f12cf9717c80f64a3547410b9a34644a7a2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: File Name / zh-cn / de-de / Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $md1Url, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $md2Url, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("D2").Value = $newZhTarget
$wsZh.Range("E2").Value = $newHandoffDatetime
$wsZh.Range("F2:G2").Clear()
$wsZh.Range("H2").Value = $newHandbackDatetime

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("D3").Value = $newZhTarget
$wsZh.Range("E3").Value = $newHandoffDatetime
$wsZh.Range("F3:G3").Clear()
$wsZh.Range("H3").Value = $newHandbackDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $md1Url, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $md1Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, $newZhTarget) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $md2Url, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $md2Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, $newZhTarget) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("D2").Value = $newDeTarget
$wsDe.Range("E2").Value = $newDeHandoffDatetime
$wsDe.Range("F2:G2").Clear()
$wsDe.Range("H2").Value = $newHandbackDatetime

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("D3").Value = $newDeTarget
$wsDe.Range("E3").Value = $newDeHandoffDatetime
$wsDe.Range("F3:G3").Clear()
$wsDe.Range("H3").Value = $newHandbackDatetime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $md1Url, [Type]::Missing, [Type]::Missing, $newFile1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $md1Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deTargetUrl, [Type]::Missing, [Type]::Missing, $newDeTarget) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $md2Url, [Type]::Missing, [Type]::Missing, $newFile2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $md2Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deTargetUrl, [Type]::Missing, [Type]::Missing, $newDeTarget) | Out-Null
